$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 29; this shifts the existing row 29 (and
# everything below it) down by one, turning the former row 29 into row 30,
# the former row 126 (last data row) into row 127, etc. — matching the
# dimension change from A1:R126 to A1:R127.
$ws.Rows.Item(29).Insert()

# Populate the newly-inserted row 29 with the new data record.
$ws.Cells.Item(29, 1).Value = 6
$ws.Cells.Item(29, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(29, 3).Value = 'Metropolitana'
$ws.Cells.Item(29, 4).Value = 44453
$ws.Cells.Item(29, 5).Value = 13
$ws.Cells.Item(29, 6).Value = 100112022
$ws.Cells.Item(29, 7).Value = 'Arveja Verde'
$ws.Cells.Item(29, 8).Value = 'Perfection'
$ws.Cells.Item(29, 9).Value = 'Primera'
$ws.Cells.Item(29, 10).Value = 220
$ws.Cells.Item(29, 11).Value = 34000
$ws.Cells.Item(29, 12).Value = 35000
$ws.Cells.Item(29, 13).Value = 34545
$ws.Cells.Item(29, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(29, 15).Value = 'Provincia de Huasco'
$ws.Cells.Item(29, 16).Value = 1382
$ws.Cells.Item(29, 17).Value = 25
$ws.Cells.Item(29, 18).Value = 'Hortaliza'
